$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

# Update D6 from "V5555" to "G5555" (adds a new shared string entry)
$ws.Range("D6").Value = "G5555"
